$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the run ", Roberto " in the "Integrantes" paragraph into two runs:
#    ", " and "Roberto " (same character formatting). We locate the text via
#    Find, then force Word to create a run boundary by toggling Bold on/off
#    over just the "Roberto " portion (net-zero formatting change, but it
#    splits the run the same way Word itself would when a user re-types part
#    of a sentence).
# ---------------------------------------------------------------------------
$rngRoberto = $d.Content
$rngRoberto.Find.Execute("Roberto ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rngRoberto.Bold = $true
$rngRoberto.Bold = $false

# ---------------------------------------------------------------------------
# 2) " Carles." -> " Carles, Jaume." (only the first occurrence, in the
#    Integrantes paragraph).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" Carles.", $true, $false, $false, $false, $false, $true, 1, $false, " Carles, Jaume.", 2)

# ---------------------------------------------------------------------------
# 3) Remove the _GoBack bookmark that currently sits at the end of the
#    "2º Sesión" paragraph (it will be re-created later at the end of the new
#    "3ª Sesión" paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4) Replace the empty "<w:ind w:left="360"/>" paragraph that follows the
#    "Empezamos con el mockup..." / second "Presentes:" paragraph with the
#    full block describing the 3rd session: the session header paragraph,
#    three bulleted list paragraphs and a final "Presentes:" paragraph.
#    We find this paragraph by locating the last empty paragraph that
#    immediately precedes the final (empty) paragraphs of the document.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "`r" -and $p.Range.ParagraphFormat.LeftIndent -eq 18) {
        $targetPara = $p
    }
}

$newXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">3ª Sesión, 04-03-2022. Inicio 19:30, Fin </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>20:30.</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Empezamos con </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Wireframe</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Nos planteamos-tras consultar al </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>profesir</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">- modificación de la idea, cambiar el </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>focus</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> de la aplicación de solamente &#8220;Pequeños negocios locales de </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>profuctos</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> relacionados con cáñamo y marihuana&#8221;, a expandir nuestro horizonte a &#8220;pequeño comercio, local y de productos/producción sostenible con el medio ambiente&#8221;.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Avanzamos el </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>wire-frame</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> desde este nuevo punto de vista.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Presentes: Sergio, </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Victor</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Gimenez</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">, </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Victor</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> Carles, Roberto </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Martinez</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$targetPara.Range.InsertXML($newXml)
